# Insert a new data row at row 700, pushing the existing rows 700:785
# down to 701:786 (dimension grows from A1:R785 to A1:R786), and fill
# in the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(700).Insert()

$ws.Range("A700").Value = 6
$ws.Range("B700").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C700").Value = "Metropolitana"
$ws.Range("D700").Value = 45124
$ws.Range("E700").Value = 13
$ws.Range("F700").Value = 100112012
$ws.Range("G700").Value = "Espinaca"
$ws.Range("H700").Value = "Sin especificar"
$ws.Range("I700").Value = "Primera"
$ws.Range("J700").Value = 410
$ws.Range("K700").Value = 7000
$ws.Range("L700").Value = 8000
$ws.Range("M700").Value = 7439
$ws.Range("N700").Value = "$/cuna 10 kilos"
$ws.Range("O700").Value = "Región Metropolitana"
$ws.Range("P700").Value = 744
$ws.Range("Q700").Value = 10
$ws.Range("R700").Value = "Hortaliza"
